$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 75: "Lista del Pueblo" / "Lista del Pueblo" (both cells styled like col B, dark/black font) ---
$ws.Range("B74").Copy()
$ws.Range("A75:B75").PasteSpecial(-4122)

# --- Row 76: "Listas Independientes" / "Listas Independientes" (both cells styled like col B) ---
$ws.Range("B74").Copy()
$ws.Range("A76:B76").PasteSpecial(-4122)

# --- Row 77: "Otros" / "Otros" (normal pattern: col A styled like existing col A, col B like col B) ---
$ws.Range("A74:B74").Copy()
$ws.Range("A77:B77").PasteSpecial(-4122)

# Fill in the values (reuses existing shared strings automatically)
$ws.Range("A75").Value = "Lista del Pueblo"
$ws.Range("B75").Value = "Lista del Pueblo"
$ws.Range("A76").Value = "Listas Independientes"
$ws.Range("B76").Value = "Listas Independientes"
$ws.Range("A77").Value = "Otros"
$ws.Range("B77").Value = "Otros"

# Clear marching-ants clipboard state from the copy operations
$excel.CutCopyMode = $false

# Move the active selection to A75, matching the saved view state
$ws.Range("A75").Select()
